$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Test - Dummy Product 7"
$ws.Range("B2").Value = "Dummy Product 7"
$ws.Range("C2").Value = "TEST - Dummy 07"
